$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: RowNumber, New DAMSLTag (col I), New DialogAct (col J)
$updates = @(
    @(2, 'b', 'Acknowledge (Backchannel)'),
    @(7, 'sd', 'Statement-non-opinion'),
    @(11, 'b', 'Acknowledge (Backchannel)'),
    @(14, 'ba', 'Appreciation'),
    @(17, '%', 'Uninterpretable'),
    @(39, 'sd', 'Statement-non-opinion'),
    @(68, 'sd', 'Statement-non-opinion'),
    @(82, 'sv', 'Statement-opinion'),
    @(96, 'sv', 'Statement-opinion'),
    @(115, 'b', 'Acknowledge (Backchannel)'),
    @(127, 'b', 'Acknowledge (Backchannel)'),
    @(158, 'sd', 'Statement-non-opinion'),
    @(161, 'ba', 'Appreciation'),
    @(169, '%', 'Uninterpretable'),
    @(174, 'sd', 'Statement-non-opinion'),
    @(175, 'sd', 'Statement-non-opinion'),
    @(180, 'sd', 'Statement-non-opinion'),
    @(197, 'sd', 'Statement-non-opinion'),
    @(201, 'aa', 'Agree/Accept'),
    @(207, 'sv', 'Statement-opinion'),
    @(213, 'aa', 'Agree/Accept'),
    @(214, 'b', 'Acknowledge (Backchannel)'),
    @(224, '%', 'Uninterpretable'),
    @(227, 'aa', 'Agree/Accept'),
    @(239, 'aa', 'Agree/Accept'),
    @(263, 'sv', 'Statement-opinion'),
    @(265, 'sv', 'Statement-opinion'),
    @(272, 'sd', 'Statement-non-opinion'),
    @(274, 'sv', 'Statement-opinion'),
    @(276, 'b', 'Acknowledge (Backchannel)'),
    @(279, '%', 'Uninterpretable'),
    @(282, 'sd', 'Statement-non-opinion'),
    @(286, 'sv', 'Statement-opinion'),
    @(287, 'aa', 'Agree/Accept'),
    @(288, 'aa', 'Agree/Accept'),
    @(292, 'ba', 'Appreciation'),
    @(294, 'sd', 'Statement-non-opinion'),
    @(296, 'b', 'Acknowledge (Backchannel)'),
    @(305, 'sv', 'Statement-opinion'),
    @(313, 'sd', 'Statement-non-opinion'),
    @(323, 'b', 'Acknowledge (Backchannel)'),
    @(325, 'b', 'Acknowledge (Backchannel)'),
    @(328, 'b', 'Acknowledge (Backchannel)'),
    @(350, 'b', 'Acknowledge (Backchannel)'),
    @(375, 'sd', 'Statement-non-opinion'),
    @(386, 'aa', 'Agree/Accept'),
    @(387, 'aa', 'Agree/Accept'),
    @(397, 'sd', 'Statement-non-opinion'),
    @(404, 'b', 'Acknowledge (Backchannel)'),
    @(414, 'aa', 'Agree/Accept'),
    @(415, 'aa', 'Agree/Accept'),
    @(426, 'aa', 'Agree/Accept'),
    @(429, 'sv', 'Statement-opinion'),
    @(439, 'sd', 'Statement-non-opinion'),
    @(440, 'sd', 'Statement-non-opinion'),
    @(451, '%', 'Uninterpretable'),
    @(456, 'aa', 'Agree/Accept'),
    @(457, 'aa', 'Agree/Accept'),
    @(458, 'b', 'Acknowledge (Backchannel)'),
    @(459, 'sv', 'Statement-opinion'),
    @(465, 'sv', 'Statement-opinion'),
    @(480, 'sv', 'Statement-opinion'),
    @(506, 'sd', 'Statement-non-opinion'),
    @(532, 'ba', 'Appreciation'),
    @(542, 'sd', 'Statement-non-opinion'),
    @(570, 'aa', 'Agree/Accept'),
    @(571, 'b', 'Acknowledge (Backchannel)'),
    @(582, 'sd', 'Statement-non-opinion'),
    @(584, 'sd', 'Statement-non-opinion'),
    @(595, 'ba', 'Appreciation'),
    @(600, 'qy', 'Yes-No-Question'),
    @(604, 'sd', 'Statement-non-opinion'),
    @(611, 'sd', 'Statement-non-opinion'),
    @(633, 'b', 'Acknowledge (Backchannel)'),
    @(650, 'sd', 'Statement-non-opinion'),
    @(676, 'sv', 'Statement-opinion'),
    @(686, 'sv', 'Statement-opinion'),
    @(689, 'sv', 'Statement-opinion'),
    @(692, 'aa', 'Agree/Accept'),
    @(699, 'sv', 'Statement-opinion'),
    @(703, 'sd', 'Statement-non-opinion'),
    @(706, 'b', 'Acknowledge (Backchannel)'),
    @(708, 'b', 'Acknowledge (Backchannel)')
)

foreach ($u in $updates) {
    $rowNum = $u[0]
    $ws.Cells.Item($rowNum, 9).Value = $u[1]
    $ws.Cells.Item($rowNum, 10).Value = $u[2]
}
